$wb = $excel.ActiveWorkbook

# --- 1. Baseline year demographics: fraction at risk of malaria 0.1 -> 0.33 ---
$wsBaseline = $wb.Worksheets.Item("Baseline year demographics")
$wsBaseline.Activate() | Out-Null
$wsBaseline.Range("C9").Value = 0.33
$wsBaseline.Range("D9").Select() | Out-Null

# --- 2. Programs cost and coverage: row 40 height + selection of entire row 44 ---
$wsCost = $wb.Worksheets.Item("Programs cost and coverage")
$wsCost.Activate() | Out-Null
$wsCost.Rows.Item(40).RowHeight = 16
$wsCost.Rows.Item(44).Select() | Out-Null

# --- 3. Programs to include: clear "x" checkmarks except Sprinkles / Treatment of MAM / Vitamin A supplementation ---
$wsInclude = $wb.Worksheets.Item("Programs to include")
$wsInclude.Activate() | Out-Null
$clearRows = @(2,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,31,32,33,34,35,36,37,38,39,41,43,45,46,47,48,49,50,51,52)
foreach ($r in $clearRows) {
    $wsInclude.Range("B$r").Value = ""
}
$wsInclude.Range("B42").Select() | Out-Null

# --- 4. Programs target population: becomes the active tab/selection ---
$wsTarget = $wb.Worksheets.Item("Programs target population")
$wsTarget.Activate() | Out-Null
$wsTarget.Range("S43").Select() | Out-Null
